# Weekly update: a new price record was reported for "Macroferia Regional de
# Talca - Piña" (Caramelo / Segunda) and inserted at row 152, pushing all
# subsequent rows (152-183) down by one (new rows 153-184). This matches how
# the source workbook is maintained (most-recent entries inserted above the
# older history for the same market/category).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the number format used by the existing date column so the newly
# inserted row keeps the same "date" styling as its neighbours.
$dateFormat = $ws.Range("D153").NumberFormat

# Insert a brand-new row above the current row 152; Excel automatically
# shifts rows 152:183 down to 153:184 (formats/styles move with them).
$ws.Rows("152:152").Insert()

# Populate the newly-inserted row with the new observation.
$ws.Range("A152").Value = 5
$ws.Range("B152").Value = "Macroferia Regional de Talca"
$ws.Range("C152").Value = "Maule"
$ws.Range("D152").Value = 44543
$ws.Range("D152").NumberFormat = $dateFormat
$ws.Range("E152").Value = 7
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100108
$ws.Range("H152").Value = "Tropicales y subtropicales"
$ws.Range("I152").Value = 100108005
$ws.Range("J152").Value = "Piña"
$ws.Range("K152").Value = "Caramelo"
$ws.Range("L152").Value = "Segunda"
$ws.Range("M152").Value = 270
$ws.Range("N152").Value = 17000
$ws.Range("O152").Value = 17000
$ws.Range("P152").Value = 17000
$ws.Range("Q152").Value = "`$/caja 14 unidades"
$ws.Range("R152").Value = "Ecuador"
$ws.Range("S152").Value = 1214
$ws.Range("T152").Value = 14
